# Updates the cryptocurrency price/volume table on Sheet1 (rows 2-51).
# For every numeric-looking "Price" (column D) value, the cell is first
# forced to Text format so Excel's COM layer stores the literal digit
# string instead of silently parsing it into a floating point number
# (which would also round/trim values such as "0.00000000360" or
# "29.992.46"). ClearFormats() afterwards drops the temporary Text
# format so the cell keeps the workbook's default (unstyled) appearance,
# exactly like the source cells before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.992.46"
$ws.Range("E2").Value = "  -2.08%  "
$ws.Range("D3").Value = "2.102.58"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -1.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "348.79"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.36%  "
$ws.Range("E6").Value = "  -0.90%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5148"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.96%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4416"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.43"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09065"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.168"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.36"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.65%  "
$ws.Range("D13").Value = "2.105.07"
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.212"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.723"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "98.85"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.74%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001147"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.003"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.66"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +6.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06658"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("E21").Value = "  -0.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.225"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.37%  "
$ws.Range("D23").Value = "30.094.47"
$ws.Range("E23").Value = "  -1.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.62"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.44%  "
$ws.Range("E25").Value = "  -1.11%  "
$ws.Range("D26").Value = "2.345.58"
$ws.Range("E26").Value = "  -0.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.96"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.585"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.14"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.32"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.171"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1059"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.654"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.229"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.960"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.139"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +4.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.13"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02565"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06782"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2276"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.51"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6828"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.318"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.67%  "
$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6511"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.44%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.16"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -5.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.278"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000360"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.14%  "
$ws.Range("E48").Value = "  -1.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.220"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "82.12"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07210"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.93%  "

